# Weekly update: a new price record (week of 2023-03-30) is inserted for
# "Arándano (blue)" at Vega Central Mapocho de Santiago, pushing the
# existing historical rows (old row 280 onward) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 280; everything at/after 280 shifts to 281+.
$ws.Rows.Item(280).Insert()

# Columns A,B,C,E,F,G,H,I,J,K,L,Q,T repeat the same "template" values as the
# row that just got pushed down to 281 (Mercado/Producto/Unidad metadata is
# identical across this product's weekly records).
$ws.Cells.Item(280, 1).Value  = $ws.Cells.Item(281, 1).Value()
$ws.Cells.Item(280, 2).Value  = $ws.Cells.Item(281, 2).Value()
$ws.Cells.Item(280, 3).Value  = $ws.Cells.Item(281, 3).Value()
$ws.Cells.Item(280, 5).Value  = $ws.Cells.Item(281, 5).Value()
$ws.Cells.Item(280, 6).Value  = $ws.Cells.Item(281, 6).Value()
$ws.Cells.Item(280, 7).Value  = $ws.Cells.Item(281, 7).Value()
$ws.Cells.Item(280, 8).Value  = $ws.Cells.Item(281, 8).Value()
$ws.Cells.Item(280, 9).Value  = $ws.Cells.Item(281, 9).Value()
$ws.Cells.Item(280, 10).Value = $ws.Cells.Item(281, 10).Value()
$ws.Cells.Item(280, 11).Value = $ws.Cells.Item(281, 11).Value()
$ws.Cells.Item(280, 12).Value = $ws.Cells.Item(281, 12).Value()
$ws.Cells.Item(280, 17).Value = $ws.Cells.Item(281, 17).Value()
$ws.Cells.Item(280, 20).Value = $ws.Cells.Item(281, 20).Value()

# New row-specific data (Fecha, Volumen, Precio mínimo/máximo/promedio,
# Origen, Precio $/Kg).
$ws.Cells.Item(280, 4).Value  = 45015
$ws.Cells.Item(280, 13).Value = 530
$ws.Cells.Item(280, 14).Value = 3800
$ws.Cells.Item(280, 15).Value = 4000
$ws.Cells.Item(280, 16).Value = 3906
$ws.Cells.Item(280, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(280, 19).Value = 1953
